$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "K" (strikes) values computed for rows 2..56 (column G), replacing the
# old "Strike#" values per the commit message ("use K instead of Strike#").
$newK = @{
    2  = 1
    3  = 2
    4  = 2
    5  = 1
    6  = 1
    7  = 0
    8  = 1
    9  = 1
    10 = 1
    11 = 0
    12 = 1
    13 = 0
    14 = 3
    15 = 0
    16 = 2
    17 = 1
    18 = 0
    19 = 1
    20 = 0
    21 = 0
    22 = 1
    23 = 0
    24 = 1
    25 = 0
    26 = 0
    27 = 2
    28 = 0
    29 = 1
    30 = 2
    31 = 2
    32 = 1
    33 = 0
    34 = 1
    35 = 3
    36 = 0
    37 = 2
    38 = 2
    39 = 3
    40 = 1
    41 = 0
    42 = 1
    43 = 0
    44 = 2
    45 = 2
    46 = 2
    47 = 3
    48 = 1
    49 = 0
    50 = 1
    51 = 0
    52 = 1
    53 = 0
    54 = 0
    55 = 2
    56 = 1
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
